$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the panel_query_time timestamp on the existing "data" sheet
$ws1.Range("F2").Value = "2021-10-05 14:19:30.780084"

# Add a new "metadata" worksheet right after the "data" sheet
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "metadata"

# Copy the header formatting (bold, bordered, centered) from the data sheet's
# header row onto the new sheet's header row (B1:G1)
$ws1.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("F1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)

# Copy the formatting of A2 (bordered/centered numeric style) for the new A2
$ws1.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row values
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row values
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "CHARGE syndrome"
$newSheet.Range("C2").Value = 745

# D2 must stay textual ("0.11"), not be auto-converted to a number, while
# keeping the default (unstyled) cell appearance
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.11"
$newSheet.Range("D2").Style = "Normal"

$newSheet.Range("E2").Value = "2019-09-25T14:37:28.081493Z"
$newSheet.Range("F2").Value = "2021-10-05 14:19:30.776407"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/745/?format=json"

# Make sure the originally-active "data" sheet stays selected/active
$ws1.Activate()
